$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of review data (row 5)
$ws.Range("A5").Value = "com.hamxa.shaynachim"
$ws.Range("B5").Value = "bitcoin guide"
$ws.Range("C5").Value = "stevewonder3001@gmail.com"
$ws.Range("D5").Value = "budoyoni@gmail.com"
$ws.Range("E5").Value = "27/5/2019 15:59"
$ws.Range("F5").Value = "I find it very helpful bitcoin guide. Try it. Great for beginners and it is free."

# Turn the two e-mail addresses into mailto hyperlinks (as in the prior rows' pattern)
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:stevewonder3001@gmail.com", "", "", "stevewonder3001@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:budoyoni@gmail.com", "", "", "budoyoni@gmail.com")

# Re-apply the same cell formatting used by the row above, restoring it after
# the hyperlink insertion (which would otherwise force its own "Hyperlink" style)
$ws.Range("A4:F4").Copy()
$ws.Range("A5:F5").PasteSpecial(-4122)

# Match the saved selection state
$ws.Range("F6").Select()
